# Updates the cryptos list with the latest price/volume(1h) figures,
# and swaps the WrappedBTC / WrappedEther rows (16 and 17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, keeping the cell unstyled
# (no explicit number format / style index), matching the source data
# which stores every Coin/Link/Price/Volume cell as a text/inline string.
function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '70.880.29'
Set-TextValue $ws.Range('E2') '  +2.12%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.566.25'
Set-TextValue $ws.Range('E3') '  +1.52%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '585.03'
Set-TextValue $ws.Range('E5') '  +1.83%  '

# Row 6
Set-TextValue $ws.Range('D6') '189.48'
Set-TextValue $ws.Range('E6') '  +1.58%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.625'
Set-TextValue $ws.Range('E7') '  +1.76%  '

# Row 8
Set-TextValue $ws.Range('D8') '3.558.88'
Set-TextValue $ws.Range('E8') '  +1.67%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.219'
Set-TextValue $ws.Range('E10') '  +15.10%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.652'
Set-TextValue $ws.Range('E11') '  +0.24%  '

# Row 12
Set-TextValue $ws.Range('D12') '54.69'
Set-TextValue $ws.Range('E12') '  +0.79%  '

# Row 13
Set-TextValue $ws.Range('D13') '0.0000318'
Set-TextValue $ws.Range('E13') '  +4.61%  '

# Row 14
Set-TextValue $ws.Range('D14') '9.53'
Set-TextValue $ws.Range('E14') '  +0.38%  '

# Row 15
Set-TextValue $ws.Range('D15') '4.132.53'
Set-TextValue $ws.Range('E15') '  +1.44%  '

# Row 16
Set-TextValue $ws.Range('B16') 'WrappedEther'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D16') '3.616.81'
Set-TextValue $ws.Range('E16') '  +2.91%  '

# Row 17
Set-TextValue $ws.Range('B17') 'WrappedBTC'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D17') '70.841.70'
Set-TextValue $ws.Range('E17') '  +2.18%  '

# Row 18
Set-TextValue $ws.Range('D18') '19.22'
Set-TextValue $ws.Range('E18') '  -0.81%  '

# Row 19
Set-TextValue $ws.Range('D19') '12.77'
Set-TextValue $ws.Range('E19') '  +3.58%  '

# Row 20
Set-TextValue $ws.Range('D20') '574.15'
Set-TextValue $ws.Range('E20') '  +5.06%  '

# Row 21
Set-TextValue $ws.Range('E21') '  +0.67%  '

# Row 22
Set-TextValue $ws.Range('D22') '1.01'
Set-TextValue $ws.Range('E22') '  -0.86%  '

# Row 23
Set-TextValue $ws.Range('D23') '17.64'
Set-TextValue $ws.Range('E23') '  -5.55%  '

# Row 24
Set-TextValue $ws.Range('D24') '4.59'
Set-TextValue $ws.Range('E24') '  +3.12%  '

# Row 25
Set-TextValue $ws.Range('D25') '4.92'
Set-TextValue $ws.Range('E25') '  -1.08%  '

# Row 26
Set-TextValue $ws.Range('D26') '94.30'
Set-TextValue $ws.Range('E26') '  -0.22%  '

# Row 27
Set-TextValue $ws.Range('D27') '11.23'
Set-TextValue $ws.Range('E27') '  -1.04%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.94'
Set-TextValue $ws.Range('E28') '  -0.51%  '

# Row 29
Set-TextValue $ws.Range('D29') '9.31'
Set-TextValue $ws.Range('E29') '  +1.56%  '

# Row 30
Set-TextValue $ws.Range('D30') '32.73'
Set-TextValue $ws.Range('E30') '  +2.47%  '

# Row 31
Set-TextValue $ws.Range('D31') '7.23'
Set-TextValue $ws.Range('E31') '  -0.79%  '

# Row 32
Set-TextValue $ws.Range('D32') '12.34'
Set-TextValue $ws.Range('E32') '  -2.53%  '

# Row 33
Set-TextValue $ws.Range('E33') '  +1.77%  '

# Row 34
Set-TextValue $ws.Range('D34') '63.70'
Set-TextValue $ws.Range('E34') '  -1.62%  '

# Row 35
Set-TextValue $ws.Range('D35') '3.81'
Set-TextValue $ws.Range('E35') '  +22.35%  '

# Row 36
Set-TextValue $ws.Range('D36') '3.31'
Set-TextValue $ws.Range('E36') '  +7.79%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.412'
Set-TextValue $ws.Range('E37') '  +1.95%  '

# Row 38
Set-TextValue $ws.Range('D38') '531.47'
Set-TextValue $ws.Range('E38') '  -4.01%  '

# Row 39
Set-TextValue $ws.Range('D39') '38.48'
Set-TextValue $ws.Range('E39') '  +0.95%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.0₃0805'
Set-TextValue $ws.Range('E40') '  +4.60%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.999'
Set-TextValue $ws.Range('E41') '  +0.04%  '

# Row 42
Set-TextValue $ws.Range('D42') '3.628.91'
Set-TextValue $ws.Range('E42') '  +9.92%  '

# Row 43
Set-TextValue $ws.Range('E43') '  +4.00%  '

# Row 44
Set-TextValue $ws.Range('D44') '3.45'
Set-TextValue $ws.Range('E44') '  +2.00%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.0470'
Set-TextValue $ws.Range('E45') '  +5.29%  '

# Row 46
Set-TextValue $ws.Range('E46') '  +0.95%  '

# Row 47
Set-TextValue $ws.Range('D47') '2.94'
Set-TextValue $ws.Range('E47') '  -1.93%  '

# Row 48
Set-TextValue $ws.Range('D48') '9.31'
Set-TextValue $ws.Range('E48') '  +4.06%  '

# Row 49
Set-TextValue $ws.Range('E49') '  +2.99%  '

# Row 51
Set-TextValue $ws.Range('E51') '  +7.36%  '
